$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New VIP time-category row
$a28 = $ws.Cells.Item(28, 1)
$b28 = $ws.Cells.Item(28, 2)
$c28 = $ws.Cells.Item(28, 3)
$d28 = $ws.Cells.Item(28, 4)

$a28.Value = "thời gian(VIP)"
$b28.Value = 250000
$c28.Value = """  """
$d28.Value = 27

# New "không được phep thay đổi tên của thời gian" note next to the "thời gian" row
$ws.Range("F27").Value = "không được phep thay đổi tên của thời gian "

# Normalise formatting on the two cells that previously held leftover blank-row
# styling so they match the rest of the table (Calibri / theme text colour,
# general alignment) instead of inheriting the old placeholder row's look.
$a28.HorizontalAlignment = 1
$a28.Font.ThemeColor = 1
$a28.Font.Name = "Calibri"

$b28.HorizontalAlignment = 1
$b28.Font.ThemeColor = 1
$b28.Font.Name = "Calibri"

$ws.Range("A28").Select()
